$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 144 (id 142) ---
# Copy formatting from the row above (143) for the styled columns (A and E)
$ws.Cells.Item(143, 1).Copy()
$ws.Cells.Item(144, 1).PasteSpecial(-4122)
$ws.Cells.Item(143, 5).Copy()
$ws.Cells.Item(144, 5).PasteSpecial(-4122)

$ws.Cells.Item(144, 1).Value = 142
$ws.Cells.Item(144, 2).Value = 7751749
$ws.Cells.Item(144, 3).Value = "India Super League"
$ws.Cells.Item(144, 4).Value = "India Super League"
$ws.Cells.Item(144, 5).Value = 45343.45833333334
$ws.Cells.Item(144, 6).Value = "FC Goa"
$ws.Cells.Item(144, 7).Value = "Northeast United"
$ws.Cells.Item(144, 11).Value = 1.4
$ws.Cells.Item(144, 12).Value = 4.2
$ws.Cells.Item(144, 13).Value = 5.75
$ws.Cells.Item(144, 14).Value = 1.4
$ws.Cells.Item(144, 15).Value = 4.2
$ws.Cells.Item(144, 16).Value = 5.75
$ws.Cells.Item(144, 17).Value = -1
$ws.Cells.Item(144, 18).Value = 1.725
$ws.Cells.Item(144, 19).Value = 2.075
$ws.Cells.Item(144, 20).Value = 2.75
$ws.Cells.Item(144, 21).Value = 1.9
$ws.Cells.Item(144, 22).Value = 1.9
$ws.Cells.Item(144, 23).Value = 0
$ws.Cells.Item(144, 24).Value = 0
$ws.Cells.Item(144, 25).Value = 0
$ws.Cells.Item(144, 26).Value = 0
$ws.Cells.Item(144, 27).Value = 0

# --- Row 145 (id 143) ---
$ws.Cells.Item(143, 1).Copy()
$ws.Cells.Item(145, 1).PasteSpecial(-4122)
$ws.Cells.Item(143, 5).Copy()
$ws.Cells.Item(145, 5).PasteSpecial(-4122)

$ws.Cells.Item(145, 1).Value = 143
$ws.Cells.Item(145, 2).Value = 7751750
$ws.Cells.Item(145, 3).Value = "India Super League"
$ws.Cells.Item(145, 4).Value = "India Super League"
$ws.Cells.Item(145, 5).Value = 45344.45833333334
$ws.Cells.Item(145, 6).Value = "Jamshedpur FC"
$ws.Cells.Item(145, 7).Value = "East Bengal Club"
$ws.Cells.Item(145, 11).Value = 2
$ws.Cells.Item(145, 12).Value = 3.3
$ws.Cells.Item(145, 13).Value = 3.4
$ws.Cells.Item(145, 14).Value = 2
$ws.Cells.Item(145, 15).Value = 3.3
$ws.Cells.Item(145, 16).Value = 3.4
$ws.Cells.Item(145, 17).Value = -0.5
$ws.Cells.Item(145, 18).Value = 2.025
$ws.Cells.Item(145, 19).Value = 1.775
$ws.Cells.Item(145, 20).Value = 2.5
$ws.Cells.Item(145, 21).Value = 1.975
$ws.Cells.Item(145, 22).Value = 1.825
$ws.Cells.Item(145, 23).Value = 0
$ws.Cells.Item(145, 24).Value = 0
$ws.Cells.Item(145, 25).Value = 0
$ws.Cells.Item(145, 26).Value = 0
$ws.Cells.Item(145, 27).Value = 0
